$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.975.17"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.882.90"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.49"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.53"
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.490"
$ws.Range("E8").Value = "  -3.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.83"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -4.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.428"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("E12").Value = "  -3.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.15"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "3.360.40"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "61.016.92"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "2.881.00"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.47"
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "423.91"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.21"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.651"
$ws.Range("E21").Value = "  -3.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.89"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.53"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.43"
$ws.Range("E24").Value = "  -4.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.04"
$ws.Range("E26").Value = "  -6.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  -3.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.07"
$ws.Range("E29").Value = "  -9.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.69"
$ws.Range("E30").Value = "  -4.58%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.63"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.104"
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("D34").Value = "0.0₃0849"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.967"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.41"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.79"
$ws.Range("E38").Value = "  -7.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.31"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.72"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.263"
$ws.Range("E43").Value = "  -6.93%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "133.00"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.648.72"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0331"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "342.46"
$ws.Range("E47").Value = "  -9.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.51"
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("E50").Value = "  -3.88%  "
$ws.Range("E51").Value = "  -3.41%  "
